$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '60.504.24'
Set-TextValue 'E2' '  +0.14%  '
Set-TextValue 'D3' '2.328.15'
Set-TextValue 'E3' '  -1.04%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '546.14'
Set-TextValue 'E5' '  -0.21%  '
Set-TextValue 'D6' '130.67'
Set-TextValue 'E6' '  -1.46%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.578'
Set-TextValue 'E8' '  -1.90%  '
Set-TextValue 'D9' '2.327.04'
Set-TextValue 'E9' '  -0.98%  '
Set-TextValue 'D10' '0.102'
Set-TextValue 'E10' '  +0.26%  '
Set-TextValue 'D11' '5.53'
Set-TextValue 'E11' '  +0.23%  '
Set-TextValue 'E12' '  -0.51%  '
Set-TextValue 'D13' '0.337'
Set-TextValue 'E13' '  +0.66%  '
Set-TextValue 'D14' '23.57'
Set-TextValue 'E14' '  -1.82%  '
Set-TextValue 'D15' '60.446.74'
Set-TextValue 'E15' '  +0.07%  '
Set-TextValue 'D16' '2.737.69'
Set-TextValue 'E16' '  -1.08%  '
Set-TextValue 'D17' '0.0000134'
Set-TextValue 'E17' '  +0.34%  '
Set-TextValue 'D18' '2.331.50'
Set-TextValue 'E18' '  -0.40%  '
Set-TextValue 'D19' '10.65'
Set-TextValue 'E19' '  -0.58%  '
Set-TextValue 'D20' '4.09'
Set-TextValue 'E20' '  -2.25%  '
Set-TextValue 'D21' '314.18'
Set-TextValue 'E21' '  -0.28%  '
Set-TextValue 'D22' '6.61'
Set-TextValue 'E22' '  -3.51%  '
Set-TextValue 'D24' '64.23'
Set-TextValue 'E24' '  +1.64%  '
Set-TextValue 'E25' '  +0.38%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.06%  '
Set-TextValue 'D27' '7.93'
Set-TextValue 'E27' '  -0.51%  '
Set-TextValue 'E28' '  +3.14%  '
Set-TextValue 'D29' '1.24'
Set-TextValue 'E29' '  +7.38%  '
Set-TextValue 'D30' '173.48'
Set-TextValue 'E30' '  +1.01%  '
Set-TextValue 'D31' '1.72'
Set-TextValue 'E31' '  -1.64%  '
Set-TextValue 'D32' '0.0₃0729'
Set-TextValue 'E32' '  -0.44%  '
Set-TextValue 'D33' '6.02'
Set-TextValue 'E33' '  +1.62%  '
Set-TextValue 'B34' 'ImmutableX'
Set-TextValue 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D34' '1.37'
Set-TextValue 'E34' '  -3.14%  '
Set-TextValue 'B35' 'PolygonEcosystemToken'
Set-TextValue 'C35' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D35' '0.383'
Set-TextValue 'E35' '  +0.18%  '
Set-TextValue 'D36' '18.04'
Set-TextValue 'E36' '  +0.01%  '
Set-TextValue 'E37' '  -0.02%  '
Set-TextValue 'D38' '1.00'
Set-TextValue 'E38' '  +0.08%  '
Set-TextValue 'D39' '4.14'
Set-TextValue 'E39' '  -0.50%  '
Set-TextValue 'D40' '321.38'
Set-TextValue 'E40' '  +1.42%  '
Set-TextValue 'D41' '38.05'
Set-TextValue 'E41' '  -0.45%  '
Set-TextValue 'D42' '1.53'
Set-TextValue 'E42' '  -0.38%  '
Set-TextValue 'D43' '137.56'
Set-TextValue 'E43' '  -3.37%  '
Set-TextValue 'D44' '3.52'
Set-TextValue 'E44' '  +1.33%  '
Set-TextValue 'D45' '0.0942'
Set-TextValue 'E45' '  -1.28%  '
Set-TextValue 'D46' '19.22'
Set-TextValue 'E46' '  -0.75%  '
Set-TextValue 'D47' '0.566'
Set-TextValue 'E47' '  +0.98%  '
Set-TextValue 'D48' '0.0498'
Set-TextValue 'E48' '  -0.25%  '
Set-TextValue 'D49' '0.0216'
Set-TextValue 'E49' '  +0.68%  '
Set-TextValue 'E50' '  +4.25%  '
Set-TextValue 'E51' '  -0.70%  '
